$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''309.81'
$ws.Range("E2").Value = '''0.96%'
$ws.Range("D3").Value = '''40.89'
$ws.Range("E3").Value = '''0.80%'
$ws.Range("D4").Value = '''5.122'
$ws.Range("E4").Value = '''0.44%'
$ws.Range("D5").Value = '''0.07676'
$ws.Range("E5").Value = '''1.18%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").Value = '''4.272'
$ws.Range("E6").Value = '''-0.05%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '''1.618'
$ws.Range("E7").Value = '''1.42%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9171'
$ws.Range("E8").Value = '''1.18%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Value = '''2.465'
$ws.Range("E9").Value = '''0.74%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1237'
$ws.Range("E10").Value = '''22.08%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1811'
$ws.Range("E11").Value = '''3.43%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09167'
$ws.Range("E12").Value = '''1.17%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.04278'
$ws.Range("E13").Value = '''2.71%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.1049'
$ws.Range("E14").Value = '''-0.51%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001247'
$ws.Range("E15").Value = '''0.73%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.005685'
$ws.Range("E16").Value = '''-3.11%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.351'
$ws.Range("E17").Value = '''0.03%'
$ws.Range("E18").Value = '''1.22%'
$ws.Range("D19").Value = '''6.891'
$ws.Range("E19").Value = '''3.89%'
$ws.Range("D20").Value = '''0.1391'
$ws.Range("E20").Value = '''2.43%'
$ws.Range("D21").Value = '''0.2729'
$ws.Range("E21").Value = '''0.00%'
$ws.Range("D22").Value = '''0.04042'
$ws.Range("E22").Value = '''-3.44%'
$ws.Range("D23").Value = '''0.001267'
$ws.Range("E23").Value = '''3.16%'
$ws.Range("D24").Value = '''0.004072'
$ws.Range("E24").Value = '''0.38%'
$ws.Range("D25").Value = '''0.0001267'
$ws.Range("E25").Value = '''-2.70%'
$ws.Range("E26").Value = '''24.43%'
$ws.Range("D38").Value = '''0.02462'
$ws.Range("E38").Value = '''2.33%'
$ws.Range("D39").Value = '''0.05258'
$ws.Range("E39").Value = '''1.89%'
$ws.Range("D40").Value = '''0.007813'
$ws.Range("E40").Value = '''0.43%'
$ws.Range("E41").Value = '''0.93%'
$ws.Range("D42").Value = '''0.006793'
$ws.Range("E42").Value = '''-3.63%'
$ws.Range("D43").Value = '''0.001838'
$ws.Range("D44").Value = '''0.008171'
$ws.Range("E44").Value = '''-3.99%'
$ws.Range("D45").Value = '''0.3093'
$ws.Range("E45").Value = '''-7.08%'
$ws.Range("D46").Value = '''0.00006829'
$ws.Range("E46").Value = '''7.27%'
$ws.Range("D47").Value = '''0.00000000748'
$ws.Range("E47").Value = '''-0.46%'
$ws.Range("D48").Value = '''0.1485'
$ws.Range("E48").Value = '''2,113.02%'
$ws.Range("D49").Value = '''0.004091'
$ws.Range("E49").Value = '''-7.24%'
$ws.Range("D50").Value = '''0.00002095'
$ws.Range("E50").Value = '''-0.46%'
$ws.Range("D51").Value = '''0.0001996'
$ws.Range("E51").Value = '''-0.46%'
